# Update countries & provincias Spain
# - Swap the Fiyi / Curazao ordering (row 199 becomes Curazao, row 200 becomes Fiyi)
#   and refresh their daily stats.
# - Refresh the "Datos actualizados..." timestamp string (cell A1).
# - Refresh a handful of per-country case numbers (Alemania row 19, Australia row 74).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Julio de 2020 a las 05:02"

# Alemania (row 19): Casos activos (D) and Recuperados (E) updated
$ws.Range("D19").Value = 187500
$ws.Range("E19").Value = 5685

# Australia (row 74): Casos totales (B), Nuevos casos (C), Casos activos (D),
# Recuperados (E), Muertes hoy (G) and Muertes (H) updated
$ws.Range("B74").Value = 11438
$ws.Range("C74").Value = 203
$ws.Range("D74").Value = 8158
$ws.Range("E74").Value = 3162
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 118

# Fiyi / Curazao swap places in the ranking (row 199 & 200) with refreshed stats.
# Row 199 now holds Curazao's updated numbers, row 200 now holds Fiyi's.
$ws.Range("A199").Value = "Curazao"
$ws.Range("B199").Value = 28
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 24
$ws.Range("E199").Value = 3
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1

$ws.Range("A200").Value = "Fiyi"
$ws.Range("B200").Value = 26
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 18
$ws.Range("E200").Value = 8
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0
